$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4273009282394885
$ws.Range("C2").Value = 0.997454690289144
$ws.Range("D2").Value = 0.5417935533659207
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor())])"
$ws.Range("G2").Value = 0.12469127785007
$ws.Range("H2").Value = 0.99
